# Generate Report for Handback
# Update the handback/handoff timestamp cells to reflect the latest run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 23:03:30"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 23:03:24"
$wsZhCn.Range("K2").Value = "2016-08-16 23:03:41"

# de-de sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-16 23:03:30"
$wsDeDe.Range("K2").Value = "2016-08-16 23:03:49"
